$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price (D) cells below contain values such as "1.007" or "0.07640" that
# Excel would otherwise reinterpret as numbers (dropping trailing zeros,
# using scientific notation, etc.) instead of keeping the original text
# formatting used by the source site. Force those cells to text first.
$priceRows = @(2, 3, 4, 5, 6, 7, 8, 9, 10, 11, 12, 13, 14, 15, 16, 18, 19, 20, 21, 22, 23, 25, 26, 27, 28, 29, 31, 32, 33, 34, 35, 37, 38, 39, 40, 41, 43, 44, 45, 46, 47, 48, 49, 50, 51)
foreach ($r in $priceRows) {
    $ws.Range("D$r").NumberFormat = "@"
}

$ws.Range("D2").Value = "26.866.55"
$ws.Range("E2").Value = "  +1.02%  "
$ws.Range("D3").Value = "1.840.65"
$ws.Range("E3").Value = "  +1.10%  "
$ws.Range("D4").Value = "1.007"
$ws.Range("E4").Value = "  -0.30%  "
$ws.Range("D5").Value = "309.14"
$ws.Range("E5").Value = "  +1.08%  "
$ws.Range("D6").Value = "1.004"
$ws.Range("E6").Value = "  -0.43%  "
$ws.Range("D7").Value = "0.4752"
$ws.Range("E7").Value = "  +1.73%  "
$ws.Range("D8").Value = "0.3671"
$ws.Range("E8").Value = "  +2.29%  "
$ws.Range("D9").Value = "0.07194"
$ws.Range("E9").Value = "  +1.02%  "
$ws.Range("D10").Value = "0.9245"
$ws.Range("E10").Value = "  +2.57%  "
$ws.Range("D11").Value = "19.60"
$ws.Range("E11").Value = "  +1.19%  "
$ws.Range("D12").Value = "0.07640"
$ws.Range("E12").Value = "  -2.04%  "
$ws.Range("D13").Value = "1.871.89"
$ws.Range("E13").Value = "  +1.84%  "
$ws.Range("D14").Value = "5.308"
$ws.Range("E14").Value = "  +1.10%  "
$ws.Range("D15").Value = "6.401"
$ws.Range("E15").Value = "  +1.07%  "
$ws.Range("D16").Value = "88.68"
$ws.Range("E16").Value = "  +1.79%  "
$ws.Range("E17").Value = "  -0.27%  "
$ws.Range("D18").Value = "0.000008638"
$ws.Range("E18").Value = "  +0.99%  "
$ws.Range("D19").Value = "1.005"
$ws.Range("E19").Value = "  -0.34%  "
$ws.Range("D20").Value = "26.897.75"
$ws.Range("E20").Value = "  +0.96%  "
$ws.Range("D21").Value = "14.52"
$ws.Range("E21").Value = "  +2.70%  "
$ws.Range("D22").Value = "5.043"
$ws.Range("E22").Value = "  +0.73%  "
$ws.Range("D23").Value = "10.64"
$ws.Range("E23").Value = "  +0.92%  "
$ws.Range("E24").Value = "  -0.97%  "
$ws.Range("D25").Value = "151.99"
$ws.Range("E25").Value = "  -0.10%  "
$ws.Range("D26").Value = "18.12"
$ws.Range("E26").Value = "  +1.35%  "
$ws.Range("D27").Value = "1.999"
$ws.Range("E27").Value = "  +1.55%  "
$ws.Range("D28").Value = "114.17"
$ws.Range("E28").Value = "  +0.63%  "
$ws.Range("D29").Value = "4.941"
$ws.Range("E29").Value = "  +3.09%  "
$ws.Range("D31").Value = "3.284"
$ws.Range("E31").Value = "  +5.10%  "
$ws.Range("D32").Value = "0.7480"
$ws.Range("E32").Value = "  +2.76%  "
$ws.Range("D33").Value = "1.173"
$ws.Range("E33").Value = "  +4.56%  "
$ws.Range("D34").Value = "4.485"
$ws.Range("E34").Value = "  +1.16%  "
$ws.Range("D35").Value = "2.746"
$ws.Range("E35").Value = "  -0.30%  "
$ws.Range("E36").Value = "  +1.20%  "
$ws.Range("D37").Value = "0.05259"
$ws.Range("E37").Value = "  +3.02%  "
$ws.Range("D38").Value = "0.01944"
$ws.Range("E38").Value = "  +1.08%  "
$ws.Range("D39").Value = "2.965"
$ws.Range("E39").Value = "  +1.67%  "
$ws.Range("D40").Value = "0.5213"
$ws.Range("E40").Value = "  +3.58%  "
$ws.Range("D41").Value = "6.935"
$ws.Range("E41").Value = "  +1.34%  "
$ws.Range("E42").Value = "  +1.29%  "
$ws.Range("D43").Value = "8.214"
$ws.Range("E43").Value = "  +2.97%  "
$ws.Range("D44").Value = "10.51"
$ws.Range("E44").Value = "  +5.85%  "
$ws.Range("D45").Value = "0.4720"
$ws.Range("E45").Value = "  +1.54%  "
$ws.Range("D46").Value = "1.004"
$ws.Range("E46").Value = "  -0.52%  "
$ws.Range("D47").Value = "101.57"
$ws.Range("E47").Value = "  +2.90%  "
$ws.Range("D48").Value = "1.603"
$ws.Range("E48").Value = "  +3.37%  "
$ws.Range("D49").Value = "65.39"
$ws.Range("E49").Value = "  +2.55%  "
$ws.Range("D50").Value = "0.06025"
$ws.Range("E50").Value = "  +0.58%  "
$ws.Range("D51").Value = "0.8847"
$ws.Range("E51").Value = "  +4.25%  "

# Restore the default (unformatted) style on the Price cells so that only
# the displayed text changes, matching the original workbook formatting.
foreach ($r in $priceRows) {
    $ws.Range("D$r").Style = "Normal"
}
